# Weekly Fruit/Vegetable price update:
# A new weekly observation is inserted as row 23 (pushing every existing
# data row at/after row 23 down by one), adding the latest price record
# for "Bruselas (repollito)" at Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 23; everything from the old row 23
# downward (through the old row 97) shifts down to rows 24-98.
$ws.Rows("23:23").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(23, 1).Value  = 9
$ws.Cells.Item(23, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(23, 3).Value  = "Metropolitana"
$ws.Cells.Item(23, 4).Value  = 45114
$ws.Cells.Item(23, 5).Value  = 13
$ws.Cells.Item(23, 6).Value  = 100112035
$ws.Cells.Item(23, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(23, 8).Value  = "Sin especificar"
$ws.Cells.Item(23, 9).Value  = "Primera"
$ws.Cells.Item(23, 10).Value = 52
$ws.Cells.Item(23, 11).Value = 19000
$ws.Cells.Item(23, 12).Value = 21000
$ws.Cells.Item(23, 13).Value = 20000
$ws.Cells.Item(23, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(23, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(23, 16).Value = 1333
$ws.Cells.Item(23, 17).Value = 15
$ws.Cells.Item(23, 18).Value = "Hortaliza"
